$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (obesity) - values updated, label unchanged
$ws.Range("B2").Value = "48.429 [-94.496 -- 180.472]"
$ws.Range("C2").Value = "51.571 [-80.472 -- 194.496]"
$ws.Range("D2").Value = "2.064 [-1.905 -- 5.623]"
$ws.Range("E2").Value = "2.198 [-2.107 -- 6.705]"
$ws.Range("F2").Value = "4.262 [0.841 -- 7.564]"
$ws.Range("G2").Value = "-0.273 [-0.816 -- 0.192]"
$ws.Range("H2").Value = "2.716 [-3.248 -- 8.504]"
$ws.Range("I2").Value = "0.728 [0.668 -- 0.832]"
$ws.Range("J2").Value = "0.779 [0.766 -- 0.798]"
$ws.Range("K2").Value = "-10.061 [-10.353 -- -9.745]"

# Row 3 (avg) - values updated, label unchanged
$ws.Range("B3").Value = "70.228 [6.284 -- 154.919]"
$ws.Range("C3").Value = "29.772 [-54.919 -- 93.716]"
$ws.Range("D3").Value = "2.555 [0.204 -- 5.010]"
$ws.Range("E3").Value = "1.083 [-0.994 -- 3.563]"
$ws.Range("F3").Value = "3.638 [0.460 -- 6.674]"
$ws.Range("G3").Value = "0.241 [-0.100 -- 0.520]"
$ws.Range("H3").Value = "1.892 [-1.503 -- 5.362]"
$ws.Range("I3").Value = "0.728 [0.668 -- 0.832]"
$ws.Range("K3").Value = "-10.061 [-10.353 -- -9.745]"

# Row 4 - label changes from "sleep disorder" to "heart disorder", values updated
$ws.Range("A4").Value = "heart disorder"
$ws.Range("B4").Value = "72.352 [-104.236 -- 151.581]"
$ws.Range("C4").Value = "27.648 [-51.581 -- 204.236]"
$ws.Range("D4").Value = "2.400 [-1.447 -- 5.871]"
$ws.Range("E4").Value = "0.917 [-1.302 -- 4.101]"
$ws.Range("F4").Value = "3.318 [0.345 -- 6.210]"
$ws.Range("G4").Value = "0.671 [-0.001 -- 1.290]"
$ws.Range("H4").Value = "2.626 [-3.032 -- 8.381]"
$ws.Range("I4").Value = "0.728 [0.668 -- 0.832]"
$ws.Range("J4").Value = "0.736 [0.716 -- 0.760]"
$ws.Range("K4").Value = "-10.061 [-10.353 -- -9.745]"

# Row 5 - label changes from "heart disorder" to "sleep disorder", values updated
$ws.Range("A5").Value = "sleep disorder"
$ws.Range("B5").Value = "95.980 [-68.777 -- 401.926]"
$ws.Range("C5").Value = "4.020 [-301.926 -- 168.777]"
$ws.Range("D5").Value = "3.200 [-0.069 -- 7.049]"
$ws.Range("E5").Value = "0.134 [-3.532 -- 4.249]"
$ws.Range("F5").Value = "3.334 [-0.195 -- 6.600]"
$ws.Range("G5").Value = "0.324 [-0.065 -- 0.692]"
$ws.Range("H5").Value = "0.335 [-6.238 -- 6.414]"
$ws.Range("I5").Value = "0.728 [0.668 -- 0.832]"
$ws.Range("J5").Value = "0.618 [0.603 -- 0.651]"
$ws.Range("K5").Value = "-10.061 [-10.353 -- -9.745]"
